$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Comments"
$ws.Range("F1").Value = "Adjustments"

$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H2").Select()
